$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.926.37"
$ws.Range("E2").Value = "  +2.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.892.45"
$ws.Range("E3").Value = "  +2.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.71"
$ws.Range("E5").Value = "  +2.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4815"
$ws.Range("E7").Value = "  +3.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2860"
$ws.Range("E8").Value = "  +5.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06548"
$ws.Range("E9").Value = "  +4.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.79"
$ws.Range("E10").Value = "  +16.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.897.24"
$ws.Range("E11").Value = "  +3.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.06"
$ws.Range("E12").Value = "  +14.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07541"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.145"
$ws.Range("E14").Value = "  +4.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6573"
$ws.Range("E15").Value = "  +6.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "300.61"
$ws.Range("E16").Value = "  +32.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.890.47"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("E18").Value = "  +6.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007563"
$ws.Range("E20").Value = "  +3.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9992"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.207"
$ws.Range("E22").Value = "  +6.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.173"
$ws.Range("E23").Value = "  +5.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.317"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.00"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.76"
$ws.Range("E26").Value = "  +11.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.964"
$ws.Range("E27").Value = "  +5.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1065"
$ws.Range("E28").Value = "  +2.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.362"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.160"
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.985"
$ws.Range("E31").Value = "  +4.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05023"
$ws.Range("E32").Value = "  +4.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.182"
$ws.Range("E33").Value = "  +3.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7286"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.714"
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01942"
$ws.Range("E36").Value = "  +3.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.719"
$ws.Range("E37").Value = "  +2.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.066"
$ws.Range("E38").Value = "  +7.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9006"
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "107.83"
$ws.Range("E40").Value = "  +3.43%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.0000"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4224"
$ws.Range("E42").Value = "  +5.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.611"
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.13"
$ws.Range("E44").Value = "  +10.33%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.389"
$ws.Range("E45").Value = "  +4.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1229"
$ws.Range("E46").Value = "  +2.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.948"
$ws.Range("E47").Value = "  +4.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.75"
$ws.Range("E48").Value = "  +5.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.398"
$ws.Range("E49").Value = "  +3.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05591"
$ws.Range("E50").Value = "  +1.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3811"
$ws.Range("E51").Value = "  +4.66%  "
